{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// 1. \"Serve {{ other_parties }}\" heading paragraph becomes a conditional:\n//    \"Serve {% if ll_name_unknown %}your landlord{% else %}{{ other_parties }}{% endif %}\"\n//    (matches the same {% if ll_name_unknown %}...{% else %}...{% endif %} pattern\n//    already used elsewhere in this template).\n// 2. The trailing empty paragraph at the very end of the document body\n//    (right before the section break) is removed.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Update the \"Serve ...\" heading paragraph -------------------------\nconst targetText = \"Serve {{ other_parties }}\";\n\nlet servePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    servePara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (servePara) {\n  const newText =\n    \"Serve {% if ll_name_unknown %}your landlord{% else %}{{ other_parties }}{% endif %}\";\n  servePara.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\n// --- 2. Remove the trailing empty paragraph at the end of the body -------\nconst count = paragraphs.items.length;\nif (count >= 2) {\n  const lastPara = paragraphs.items[count - 1];\n  const secondLastPara = paragraphs.items[count - 2];\n  if (lastPara.text === \"\") {\n    const startRange = secondLastPara.getRange(\"End\");\n    const endRange = lastPara.getRange(\"End\");\n    const spanRange = startRange.expandTo(endRange);\n    spanRange.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# 1. \"Serve {{ other_parties }}\" heading paragraph becomes a conditional:\n#    \"Serve {% if ll_name_unknown %}your landlord{% else %}{{ other_parties }}{% endif %}\"\n#    (matches the same {% if ll_name_unknown %}...{% else %}...{% endif %} pattern\n#    already used elsewhere in this template).\n# 2. The trailing empty paragraph at the very end of the document body\n#    (right before the section break) is removed.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the \"Serve ...\" heading paragraph --------------------------\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Serve {{ other_parties }}\")\nif ($found) {\n    $findRange.Text = \"Serve {% if ll_name_unknown %}your landlord{% else %}{{ other_parties }}{% endif %}\"\n}\n\n# --- 2. Remove the trailing empty paragraph at the end of the document ----\n$count = $d.Paragraphs.Count\nif ($count -ge 2) {\n    $lastPara = $d.Paragraphs.Item($count)\n    $secondLastPara = $d.Paragraphs.Item($count - 1)\n    if ($lastPara.Range.Text -eq [char]13) {\n        # Delete the paragraph mark that separates the second-to-last\n        # paragraph from the trailing empty one, merging them and removing\n        # the empty paragraph while keeping the previous paragraph's text.\n        $mark = $d.Range($secondLastPara.Range.End - 1, $secondLastPara.Range.End)\n        $mark.Delete()\n    }\n}\n"}
